# redmine #9271 - cal sheet update for GA05MOAS-GL496 (Omaha_Cal_Info)
#
# Real content edits captured by this script:
#   1) Moorings!J2 (Cruise Number) "AT26-30" -> "AT-26-30": a missing hyphen
#      after the ship code is inserted; the newly-inserted "-" is flagged in
#      blue so the correction stands out (rich text run).
#   2) Asset_Cal_Info!F4 (CC_angular_resolution coefficient) value corrected
#      from 1.13 to 1.096, flagged in blue to mark the corrected value.
#   3) View-state touch-ups left behind by the resave (active cell / zoom).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Moorings"
# ---------------------------------------------------------------------
$moorings = $wb.Worksheets.Item("Moorings")
$moorings.Select()

$cruiseCell = $moorings.Cells.Item(2, 10)   # J2, "Cruise Number"
$cruiseCell.Value = "AT-26-30"

# Give the trailing part ("26-30") the same explicit Arial 12 run the
# source workbook carries, inheriting automatic (black) colour.
$tail = $cruiseCell.Characters(4, 5)
$tail.Font.Name = "Arial"
$tail.Font.Size = 12
$tail.Font.ColorIndex = -4105

# Colour just the newly-inserted hyphen blue to flag the correction.
$hyphen = $cruiseCell.Characters(3, 1)
$hyphen.Font.Color = 255 * 65536

# Restore default (100%) zoom and move the active selection, matching the
# post-edit view state.
$excel.ActiveWindow.Zoom = 100
$moorings.Range("D11").Select()

# ---------------------------------------------------------------------
# Sheet "Asset_Cal_Info"
# ---------------------------------------------------------------------
$assetCal = $wb.Worksheets.Item("Asset_Cal_Info")
$assetCal.Select()

$coeffCell = $assetCal.Cells.Item(4, 6)   # F4, CC_angular_resolution value
$coeffCell.Value = 1.096
$coeffCell.Font.Color = 255 * 65536   # flag the corrected value in blue

$assetCal.Range("F4").Select()

$moorings.Select()
